$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new leading data row (winter 2007 forecast) was added to the series, which
# pushes all the existing data rows down by one. Insert a blank row at row 2
# so the old row 2 becomes row 3, old row 3 becomes row 4, etc.
$ws.Rows.Item(2).Insert()

# Insert() copies formatting down from the row above (the bold header row);
# strip that back out so the new row starts from the default/plain style.
$ws.Rows.Item(2).ClearFormats()

# Column A holds dates and uses a dedicated date-formatted style throughout
# the table (same style as A3:A19). Copy just that formatting onto A2.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-simulated rt_data: refresh every row's values (dates/years are carried
# over/shifted by one row; the y_0_forecast (C) and y_1_forecast (E) columns
# were recomputed for every row as part of the evaluation bugfix).
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = -0.7015558851707349
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 0.3143490788445336

$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = 0.3590181115727287
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = 0.6970543652217165

$ws.Range("A4").Value = 40130
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = -0.01587181126743165
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = -0.1225239755399454

$ws.Range("A5").Value = 40494
$ws.Range("B5").Value = 2010
$ws.Range("C5").Value = -0.02256889165885845
$ws.Range("D5").Value = 2011
$ws.Range("E5").Value = -0.02753509623226735

$ws.Range("A6").Value = 40862
$ws.Range("B6").Value = 2011
$ws.Range("C6").Value = 0.09611428386597787
$ws.Range("D6").Value = 2012
$ws.Range("E6").Value = -0.2256894584805158

$ws.Range("A7").Value = 41228
$ws.Range("B7").Value = 2012
$ws.Range("C7").Value = -0.1827723404408288
$ws.Range("D7").Value = 2013
$ws.Range("E7").Value = -0.2251688766575

$ws.Range("A8").Value = 41592
$ws.Range("B8").Value = 2013
$ws.Range("C8").Value = -0.001350220946483294
$ws.Range("D8").Value = 2014
$ws.Range("E8").Value = 0.1494097328869959

$ws.Range("A9").Value = 41957
$ws.Range("B9").Value = 2014
$ws.Range("C9").Value = -0.075754880139145
$ws.Range("D9").Value = 2015
$ws.Range("E9").Value = -0.5497151367044428

$ws.Range("A10").Value = 42321
$ws.Range("B10").Value = 2015
$ws.Range("C10").Value = -0.5761528471665445
$ws.Range("D10").Value = 2016
$ws.Range("E10").Value = -0.1255150964614482

$ws.Range("A11").Value = 42689
$ws.Range("B11").Value = 2016
$ws.Range("C11").Value = -0.2011999787958185
$ws.Range("D11").Value = 2017
$ws.Range("E11").Value = -0.3746351385105373

$ws.Range("A12").Value = 43053
$ws.Range("B12").Value = 2017
$ws.Range("C12").Value = 0.1213692818849532
$ws.Range("D12").Value = 2018
$ws.Range("E12").Value = 0.1731436979489631

$ws.Range("A13").Value = 43418
$ws.Range("B13").Value = 2018
$ws.Range("C13").Value = 0.1493219406571766
$ws.Range("D13").Value = 2019
$ws.Range("E13").Value = 0.9013851022877439

$ws.Range("A14").Value = 43783
$ws.Range("B14").Value = 2019
$ws.Range("C14").Value = -0.4278219446121612
$ws.Range("D14").Value = 2020
$ws.Range("E14").Value = -0.8258413506386342

$ws.Range("A15").Value = 44159
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = -1.026566979837418
$ws.Range("D15").Value = 2021
$ws.Range("E15").Value = -2.816143384276215

$ws.Range("A16").Value = 44525
$ws.Range("B16").Value = 2021
$ws.Range("C16").Value = 0.3179894933462268
$ws.Range("D16").Value = 2022
$ws.Range("E16").Value = 0.07011423530434158

$ws.Range("A17").Value = 44890
$ws.Range("B17").Value = 2022
$ws.Range("C17").Value = 0.463604920919658
$ws.Range("D17").Value = 2023
$ws.Range("E17").Value = -0.3203420516749933

$ws.Range("A18").Value = 45254
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = 0.621639092134818
$ws.Range("D18").Value = 2024
$ws.Range("E18").Value = 0.009137938461889483

$ws.Range("A19").Value = 45618
$ws.Range("B19").Value = 2024
$ws.Range("C19").Value = -0.6768900623516982
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = 0.9453792747973422
